# "anyadido tipo y cuota recargo equivalencia"
# (added equivalence-surcharge type & fee) - correct invoice A13069 -> A13095
# and rework the amounts on row 2: fix Base1/Cuota1, drop the old
# TipoRE1/CuotaRE1 (equivalence surcharge) pair from columns H/I and
# populate the Base2/Tipo2/Cuota2 and Base3/Tipo3/Cuota3 blocks instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Correct the invoice number (NumFactura) in A2
$ws.Range("A2").Value = "A13095"

# Update Base1 (E2) / Cuota1 (G2) values
$ws.Range("E2").Value = 2.13
$ws.Range("G2").Value = 0.45

# Remove the old TipoRE1 (H2) / CuotaRE1 (I2) values entirely
$ws.Range("H2:I2").Clear()

# Populate the new Base2 (J2) / Tipo2 (K2) / Cuota2 (L2) block
$ws.Range("J2").Value = 113.07
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 11.31

# Populate the new Base3 (O2) / Tipo3 (P2) / Cuota3 (Q2) block
$ws.Range("O2").Value = 100
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 5
